# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new row 2 for the 2022-Q3 totals,
#    pushing the existing 2022-Q2 totals row down to row 3.
# 2. Duplicate the existing "2022-Q2" holdings sheet so the old data is kept
#    intact on its own tab (it becomes the new, third, "2022-Q2" sheet).
# 3. Replace the contents of the original "2022-Q2" sheet (which stays in
#    the 2nd tab position) with the new 2022-Q3 fund-holdings table and
#    rename that tab to "2022-Q3".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet — add the 2022-Q3 summary row above the existing 2022-Q2 row
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Re-use the formatting of the row that just got pushed down for the new row
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "'2022-Q3"
$summary.Cells.Item(2, 3).Value = 14
$summary.Cells.Item(2, 4).Value = 5.78

$summary.Cells.Item(3, 1).Value = 1

# ---------------------------------------------------------------------------
# 2. Preserve the old 2022-Q2 fund-holdings sheet by duplicating it as a new
#    tab before we overwrite the original with 2022-Q3 data.
# ---------------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Copy($null, $oldQ2)
$preservedQ2 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 3. Turn the original (2nd tab) sheet into the 2022-Q3 table.
# ---------------------------------------------------------------------------
$q3 = $oldQ2
$q3.Name = "2022-Q3"
$preservedQ2.Name = "2022-Q2"

# Wipe the old values but keep the cell formatting (borders/bold header etc.)
$q3.Cells.ClearContents()

# The new table has one extra data row (14 funds vs. 13 before) — extend the
# bold/bordered "index" column formatting down into that new row 15.
$q3.Cells.Item(14, 1).Copy()
$q3.Cells.Item(15, 1).PasteSpecial(-4122)

$q3.Cells.Item(1, 2).Value = "'基金代码"
$q3.Cells.Item(1, 3).Value = "'基金名称"
$q3.Cells.Item(1, 4).Value = "'基金规模"
$q3.Cells.Item(1, 5).Value = "'股票总仓位"
$q3.Cells.Item(1, 6).Value = "'仓位占比"
$q3.Cells.Item(1, 7).Value = "'持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "'仓位排名"
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "'159941"
$q3.Cells.Item(2, 3).Value = "'广发纳斯达克100ETF（QDII）"
$q3.Cells.Item(2, 4).Value = "'106.15"
$q3.Cells.Item(2, 5).Value = "'91.14"
$q3.Cells.Item(2, 6).Value = "'2.03"
$q3.Cells.Item(2, 7).Value = "'2.1548"
$q3.Cells.Item(2, 8).Value = 9
$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "'513100"
$q3.Cells.Item(3, 3).Value = "'国泰纳斯达克100（QDII-ETF）"
$q3.Cells.Item(3, 4).Value = "'46.54"
$q3.Cells.Item(3, 5).Value = "'91.35"
$q3.Cells.Item(3, 6).Value = "'2.03"
$q3.Cells.Item(3, 7).Value = "'0.9448"
$q3.Cells.Item(3, 8).Value = 9
$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "'040047"
$q3.Cells.Item(4, 3).Value = "'华安纳斯达克100指数（QDII）美元现钞A"
$q3.Cells.Item(4, 4).Value = "'24.52"
$q3.Cells.Item(4, 5).Value = "'92.09"
$q3.Cells.Item(4, 6).Value = "'2.06"
$q3.Cells.Item(4, 7).Value = "'0.5051"
$q3.Cells.Item(4, 8).Value = 9
$q3.Cells.Item(5, 1).Value = 3
$q3.Cells.Item(5, 2).Value = "'040048"
$q3.Cells.Item(5, 3).Value = "'华安纳斯达克100指数（QDII）美元现汇A"
$q3.Cells.Item(5, 4).Value = "'24.52"
$q3.Cells.Item(5, 5).Value = "'92.09"
$q3.Cells.Item(5, 6).Value = "'2.06"
$q3.Cells.Item(5, 7).Value = "'0.5051"
$q3.Cells.Item(5, 8).Value = 9
$q3.Cells.Item(6, 1).Value = 4
$q3.Cells.Item(6, 2).Value = "'040046"
$q3.Cells.Item(6, 3).Value = "'华安纳斯达克100指数（QDII）人民币A"
$q3.Cells.Item(6, 4).Value = "'22.21"
$q3.Cells.Item(6, 5).Value = "'92.09"
$q3.Cells.Item(6, 6).Value = "'2.06"
$q3.Cells.Item(6, 7).Value = "'0.4575"
$q3.Cells.Item(6, 8).Value = 9
$q3.Cells.Item(7, 1).Value = 5
$q3.Cells.Item(7, 2).Value = "'160213"
$q3.Cells.Item(7, 3).Value = "'国泰纳斯达克100指数（QDII）"
$q3.Cells.Item(7, 4).Value = "'15.14"
$q3.Cells.Item(7, 5).Value = "'85.81"
$q3.Cells.Item(7, 6).Value = "'1.96"
$q3.Cells.Item(7, 7).Value = "'0.2967"
$q3.Cells.Item(7, 8).Value = 9
$q3.Cells.Item(8, 1).Value = 6
$q3.Cells.Item(8, 2).Value = "'000834"
$q3.Cells.Item(8, 3).Value = "'大成纳斯达克100指数（QDII）"
$q3.Cells.Item(8, 4).Value = "'14.15"
$q3.Cells.Item(8, 5).Value = "'85.22"
$q3.Cells.Item(8, 6).Value = "'1.90"
$q3.Cells.Item(8, 7).Value = "'0.2688"
$q3.Cells.Item(8, 8).Value = 9
$q3.Cells.Item(9, 1).Value = 7
$q3.Cells.Item(9, 2).Value = "'513300"
$q3.Cells.Item(9, 3).Value = "'华夏纳斯达克100ETF（QDII）"
$q3.Cells.Item(9, 4).Value = "'11.08"
$q3.Cells.Item(9, 5).Value = "'97.32"
$q3.Cells.Item(9, 6).Value = "'2.18"
$q3.Cells.Item(9, 7).Value = "'0.2415"
$q3.Cells.Item(9, 8).Value = 3
$q3.Cells.Item(10, 1).Value = 8
$q3.Cells.Item(10, 2).Value = "'003722"
$q3.Cells.Item(10, 3).Value = "'易方达纳斯达克100指数美元（QDII-LOF）A"
$q3.Cells.Item(10, 4).Value = "'7.72"
$q3.Cells.Item(10, 5).Value = "'90.67"
$q3.Cells.Item(10, 6).Value = "'2.05"
$q3.Cells.Item(10, 7).Value = "'0.1583"
$q3.Cells.Item(10, 8).Value = 9
$q3.Cells.Item(11, 1).Value = 9
$q3.Cells.Item(11, 2).Value = "'161130"
$q3.Cells.Item(11, 3).Value = "'易方达纳斯达克100指数人民币（QDII-LOF）"
$q3.Cells.Item(11, 4).Value = "'7.72"
$q3.Cells.Item(11, 5).Value = "'90.67"
$q3.Cells.Item(11, 6).Value = "'2.05"
$q3.Cells.Item(11, 7).Value = "'0.1583"
$q3.Cells.Item(11, 8).Value = 9
$q3.Cells.Item(12, 1).Value = 10
$q3.Cells.Item(12, 2).Value = "'014978"
$q3.Cells.Item(12, 3).Value = "'华安纳斯达克100指数（QDII）人民币C"
$q3.Cells.Item(12, 4).Value = "'2.31"
$q3.Cells.Item(12, 5).Value = "'92.09"
$q3.Cells.Item(12, 6).Value = "'2.06"
$q3.Cells.Item(12, 7).Value = "'0.0476"
$q3.Cells.Item(12, 8).Value = 9
$q3.Cells.Item(13, 1).Value = 11
$q3.Cells.Item(13, 2).Value = "'159632"
$q3.Cells.Item(13, 3).Value = "'华安纳斯达克100ETF（QDII）"
$q3.Cells.Item(13, 4).Value = "'1.51"
$q3.Cells.Item(13, 5).Value = "'89.05"
$q3.Cells.Item(13, 6).Value = "'2.00"
$q3.Cells.Item(13, 7).Value = "'0.0302"
$q3.Cells.Item(13, 8).Value = 9
$q3.Cells.Item(14, 1).Value = 12
$q3.Cells.Item(14, 2).Value = "'012871"
$q3.Cells.Item(14, 3).Value = "'易方达纳斯达克100指数美元（QDII-LOF）C"
$q3.Cells.Item(14, 4).Value = "'0.18"
$q3.Cells.Item(14, 5).Value = "'90.67"
$q3.Cells.Item(14, 6).Value = "'2.05"
$q3.Cells.Item(14, 7).Value = "'0.0037"
$q3.Cells.Item(14, 8).Value = 9
$q3.Cells.Item(15, 1).Value = 13
$q3.Cells.Item(15, 2).Value = "'012870"
$q3.Cells.Item(15, 3).Value = "'易方达纳斯达克100指数人民币（QDII-LOF）C"
$q3.Cells.Item(15, 4).Value = "'0.18"
$q3.Cells.Item(15, 5).Value = "'90.67"
$q3.Cells.Item(15, 6).Value = "'2.05"
$q3.Cells.Item(15, 7).Value = "'0.0037"
$q3.Cells.Item(15, 8).Value = 9

